$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A16").Value = "LFU"
$ws.Range("B16").Value = 31101562
$ws.Range("C16").Value = "TournamentBP"
